$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '25.912.71'
$ws.Range("E2").Value = "'" + '  -1.39%  '
$ws.Range("D3").Value = "'" + '1.637.26'
$ws.Range("E3").Value = "'" + '  -0.68%  '
$ws.Range("E4").Value = "'" + '  +0.25%  '
$ws.Range("D5").Value = "'" + '214.85'
$ws.Range("E5").Value = "'" + '  -1.09%  '
$ws.Range("D6").Value = "'" + '0.505'
$ws.Range("E6").Value = "'" + '  -0.22%  '
$ws.Range("E7").Value = "'" + '  +0.27%  '
$ws.Range("E8").Value = "'" + '  -1.10%  '
$ws.Range("D9").Value = "'" + '0.0639'
$ws.Range("E9").Value = "'" + '  +0.09%  '
$ws.Range("D10").Value = "'" + '19.59'
$ws.Range("E10").Value = "'" + '  -2.42%  '
$ws.Range("D11").Value = "'" + '0.0793'
$ws.Range("E11").Value = "'" + '  -0.25%  '
$ws.Range("B12").Value = "'" + 'WrappedEther'
$ws.Range("C12").Value = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'" + '1.772.58'
$ws.Range("E12").Value = "'" + '  +8.37%  '
$ws.Range("B13").Value = "'" + 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = "'" + 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = "'" + '1.864.35'
$ws.Range("E13").Value = "'" + '  -0.66%  '
$ws.Range("E14").Value = "'" + '  -1.16%  '
$ws.Range("E15").Value = "'" + '  -2.05%  '
$ws.Range("D16").Value = "'" + '0.0₃0765'
$ws.Range("E16").Value = "'" + '  -0.49%  '
$ws.Range("D17").Value = "'" + '62.75'
$ws.Range("E17").Value = "'" + '  -1.65%  '
$ws.Range("D18").Value = "'" + '25.922.96'
$ws.Range("E18").Value = "'" + '  -1.29%  '
$ws.Range("E19").Value = "'" + '  +0.32%  '
$ws.Range("D20").Value = "'" + '193.00'
$ws.Range("E20").Value = "'" + '  -1.98%  '
$ws.Range("E21").Value = "'" + '  -1.91%  '
$ws.Range("E22").Value = "'" + '  -1.45%  '
$ws.Range("E23").Value = "'" + '  -1.31%  '
$ws.Range("D24").Value = "'" + '144.03'
$ws.Range("E24").Value = "'" + '  +0.35%  '
$ws.Range("E25").Value = "'" + '  +0.34%  '
$ws.Range("E27").Value = "'" + '  +0.90%  '
$ws.Range("D28").Value = "'" + '6.84'
$ws.Range("E28").Value = "'" + '  -1.76%  '
$ws.Range("E29").Value = "'" + '  -1.15%  '
$ws.Range("E30").Value = "'" + '  -1.02%  '
$ws.Range("D31").Value = "'" + '0.0501'
$ws.Range("E31").Value = "'" + '  -0.18%  '
$ws.Range("E32").Value = "'" + '  -1.82%  '
$ws.Range("E33").Value = "'" + '  -0.99%  '
$ws.Range("D34").Value = "'" + '1.54'
$ws.Range("E34").Value = "'" + '  -4.12%  '
$ws.Range("E35").Value = "'" + '  +1.05%  '
$ws.Range("E36").Value = "'" + '  -1.68%  '
$ws.Range("D37").Value = "'" + '1.137.64'
$ws.Range("E37").Value = "'" + '  -0.17%  '
$ws.Range("E38").Value = "'" + '  -2.48%  '
$ws.Range("E39").Value = "'" + '  -1.67%  '
$ws.Range("E40").Value = "'" + '  -0.50%  '
$ws.Range("E41").Value = "'" + '  +0.23%  '
$ws.Range("B42").Value = "'" + 'FraxShare'
$ws.Range("C42").Value = "'" + 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'" + '5.47'
$ws.Range("E42").Value = "'" + '  -3.37%  '
$ws.Range("D43").Value = "'" + '99.36'
$ws.Range("E43").Value = "'" + '  -1.22%  '
$ws.Range("B44").Value = "'" + 'TrustWalletToken'
$ws.Range("C44").Value = "'" + 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = "'" + '0.799'
$ws.Range("E44").Value = "'" + '  -0.32%  '
$ws.Range("D45").Value = "'" + '1.774.35'
$ws.Range("E45").Value = "'" + '  -0.64%  '
$ws.Range("E46").Value = "'" + '  +2.68%  '
$ws.Range("D47").Value = "'" + '56.53'
$ws.Range("E47").Value = "'" + '  +0.20%  '
$ws.Range("E48").Value = "'" + '  +2.53%  '
$ws.Range("E49").Value = "'" + '  -1.03%  '
$ws.Range("D50").Value = "'" + '7.68'
$ws.Range("E50").Value = "'" + '  -0.29%  '
$ws.Range("E51").Value = "'" + '  -0.70%  '
